# Before using the code edit the config file
#
# The sheet originally has six "status" columns (E:J), each header cell
# ("status") carrying its own green fill, and five data rows below it
# filled with "Pass". This change extends that block with eleven more
# identical "status" columns, K through U, reusing the same two shared
# strings ("status" / "Pass") and giving every new header cell its own
# (new) fill so the style table grows just like the original did when
# each column was added one at a time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$destCols = @("K","L","M","N","O","P","Q","R","S","T","U")

# Matches the width already used by the existing status columns (E:J)
# (their stored <col> width is 6.3515625 units; 5.5 "characters" is the
# closest the ColumnWidth setter's pixel grid can land on that value).
$statusWidth = 5.5

$colorSeed = 32768  # RGB 0x008000 ("dark green"), same hue as the existing status fills

for ($i = 0; $i -lt $destCols.Length; $i++) {
    $col = $destCols[$i]

    $header = $ws.Range($col + "1")
    $header.Value = "status"
    $header.ColumnWidth = $statusWidth
    # Give every new header its own distinct fill entry (a new style slot),
    # same as each pre-existing status column E:J already has.
    $header.Interior.Color = $colorSeed + $i

    for ($r = 2; $r -le 6; $r++) {
        $ws.Range($col + $r).Value = "Pass"
    }
}
